# Remove the "Guessed License URL Text" / $guessedLicenseContent$ column
# (column P) from the ComponentsAndLicenses sheet. Deleting the whole
# column shifts every column to its right one position to the left and
# drops the now-unused shared strings, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComponentsAndLicenses")

$ws.Columns("P:P").Delete()

$ws.Activate()
$ws.Range("W8").Select()
